# Auto-applied edit: "added 4wk low sales check"
# Forecast Comparison sheet: MyForecast/Trend/Inventory Coverage/Stockout Risk/
# Reorder Urgency/Seasonality Index updated for the low-sales-volume recalculation.
# Summary sheet: forecast totals recomputed to match.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row 2
$wsForecast.Range("D2").Value = 1
$wsForecast.Range("G2").Value = "Low Volume Season"
$wsForecast.Range("H2").Value = 20.56
$wsForecast.Range("L2").Value = 0.9

# Row 3
$wsForecast.Range("D3").Value = 1
$wsForecast.Range("G3").Value = "Low Volume Season"
$wsForecast.Range("H3").Value = 19.56
$wsForecast.Range("L3").Value = 1.1

# Row 4
$wsForecast.Range("D4").Value = 1
$wsForecast.Range("G4").Value = "Low Volume Season"
$wsForecast.Range("H4").Value = 18.56
$wsForecast.Range("L4").Value = 0.91

# Row 5
$wsForecast.Range("D5").Value = 1
$wsForecast.Range("G5").Value = "Low Volume Season"
$wsForecast.Range("H5").Value = 17.56
$wsForecast.Range("L5").Value = 0.84

# Row 6
$wsForecast.Range("D6").Value = 1
$wsForecast.Range("G6").Value = "Low Volume Season"
$wsForecast.Range("H6").Value = 16.56
$wsForecast.Range("L6").Value = 0.99

# Row 7
$wsForecast.Range("D7").Value = 1
$wsForecast.Range("G7").Value = "Low Volume Season"
$wsForecast.Range("H7").Value = 15.56
$wsForecast.Range("L7").Value = 0.88

# Row 8
$wsForecast.Range("D8").Value = 44
$wsForecast.Range("G8").Value = "Low Volume Season"
$wsForecast.Range("H8").Value = 0.6
$wsForecast.Range("J8").Value = "Urgent"
$wsForecast.Range("L8").Value = 1

# Row 9
$wsForecast.Range("D9").Value = 61
$wsForecast.Range("G9").Value = "Low Volume Season"
$wsForecast.Range("H9").Value = 0
$wsForecast.Range("I9").Value = "High"
$wsForecast.Range("J9").Value = "Urgent"
$wsForecast.Range("L9").Value = 0.83

# Row 10
$wsForecast.Range("D10").Value = 21
$wsForecast.Range("G10").Value = "Low Volume Season"
$wsForecast.Range("H10").Value = 0
$wsForecast.Range("I10").Value = "High"
$wsForecast.Range("J10").Value = "Urgent"
$wsForecast.Range("L10").Value = 1.14

# Row 11
$wsForecast.Range("D11").Value = 1
$wsForecast.Range("G11").Value = "Low Volume Season"
$wsForecast.Range("H11").Value = 0
$wsForecast.Range("I11").Value = "High"
$wsForecast.Range("J11").Value = "Urgent"
$wsForecast.Range("L11").Value = 1.11

# Row 12
$wsForecast.Range("D12").Value = 1
$wsForecast.Range("G12").Value = "Low Volume Season"
$wsForecast.Range("H12").Value = 0
$wsForecast.Range("I12").Value = "High"
$wsForecast.Range("L12").Value = 1.14

# Row 13
$wsForecast.Range("D13").Value = 1
$wsForecast.Range("G13").Value = "Low Volume Season"
$wsForecast.Range("L13").Value = 0.84

# Row 14
$wsForecast.Range("D14").Value = 35
$wsForecast.Range("G14").Value = "Low Volume Season"
$wsForecast.Range("L14").Value = 0.92

# Row 15
$wsForecast.Range("D15").Value = 66
$wsForecast.Range("G15").Value = "Low Volume Season"
$wsForecast.Range("L15").Value = 0.9399999999999999

# Row 16
$wsForecast.Range("D16").Value = 26
$wsForecast.Range("G16").Value = "Low Volume Season"
$wsForecast.Range("L16").Value = 0.83

# Row 17
$wsForecast.Range("D17").Value = 1
$wsForecast.Range("G17").Value = "Low Volume Season"
$wsForecast.Range("L17").Value = 1.12

# Summary sheet forecast totals (stored as text, matching the sheet's existing convention)
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "273"
$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "116"
$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "7"
$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "66"
$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "2"
